$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename Sheet1 -> "Moment Club"
$ws.Name = "Moment Club"

# Add new "GalleryImageList" header column (AU) and per-row gallery JSON values.
$galleryJson = '["D:/Web Automations/cypress/fixtures/Files/Host Images/profile.jpg","D:/Web Automations/cypress/fixtures/Files/Host Images/host1.jpg","D:/Web Automations/cypress/fixtures/Files/Host Images/host2.jpeg"]'

$ws.Range("AU1").Value = "GalleryImageList"
$ws.Range("AU2").Value = $galleryJson
$ws.Range("AU3").Value = $galleryJson
$ws.Range("AU4").Value = $galleryJson

# Add host image path values for row 2 (host profile photo paths)
$ws.Range("C2").Value = "D:/Web Automations/cypress/fixtures/Files/Host Images/host1.jpg"
$ws.Range("D2").Value = "D:/Web Automations/cypress/fixtures/Files/Host Images/host3.png\"

# Update the active cell selection shown in the saved workbook view
$ws.Range("AI17").Select()
